$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A123").Value = 45935
$ws.Range("B123").Value = "#1 - Vajtswv Thov Koj Foom Koobhmoov"
$ws.Range("C123").Value = "Public - No reporting needed"
$ws.Range("D123").Value = "HBNA Songbook"
$ws.Range("E123").Value = "Guest Speaker Week (Unknown Topic)"

Write-Host "done"
